# Apply the cryptos-list refresh described in the commit diff.
# For each changed row we update the "Price" (column D) and/or "Volume(1h)"
# (column E) cell to the new scraped text.
#
# NOTE: column D holds plain-text price strings (e.g. "70.734.59", "1.00",
# "0.155") -- NOT real numbers (some even use "." as a thousands separator).
# Assigning such look-alike-numeric strings straight to Range.Value makes Excel
# silently reinterpret them as numbers (losing the original formatting/trailing
# zeros). To keep them as genuine text -- matching the workbook before the edit --
# we write a formula that evaluates to the literal string and then immediately
# Copy / PasteSpecial(values) it back onto itself; this freezes the formula result
# as a static text cell with no left-over formula and no style changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.734.59"
$ws.Range("E2").Value = "  +0.64%  "
$ws.Range("D3").Value = "3.531.08"
$ws.Range("E3").Value = "  -0.35%  "
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Formula = "=""620.79"""
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  +4.09%  "
$ws.Range("D6").Formula = "=""172.62"""
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  +0.55%  "
$ws.Range("D7").Value = "3.526.14"
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("E8").Value = "  -0.41%  "
$ws.Range("D9").Formula = "=""1.00"""
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("E10").Value = "  +1.83%  "
$ws.Range("D11").Formula = "=""7.21"""
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = "  -2.88%  "
$ws.Range("E12").Value = "  +0.50%  "
$ws.Range("D13").Formula = "=""46.19"""
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = "  -0.24%  "
$ws.Range("E14").Value = "  -0.20%  "
$ws.Range("D15").Value = "4.099.01"
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("D16").Formula = "=""8.43"""
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = "  +1.61%  "
$ws.Range("D17").Formula = "=""607.06"""
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value = "  -0.34%  "
$ws.Range("D18").Value = "3.540.65"
$ws.Range("E18").Value = "  -0.63%  "
$ws.Range("D19").Value = "70.840.96"
$ws.Range("E19").Value = "  +0.62%  "
$ws.Range("E20").Value = "  +1.81%  "
$ws.Range("D21").Formula = "=""17.67"""
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  +2.37%  "
$ws.Range("E22").Value = "  +0.54%  "
$ws.Range("D23").Formula = "=""9.07"""
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  -0.94%  "
$ws.Range("D24").Formula = "=""15.65"""
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("D25").Formula = "=""97.75"""
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  +1.39%  "
$ws.Range("E26").Value = "  -0.31%  "
$ws.Range("E27").Value = "  +0.00%  "
$ws.Range("E28").Value = "  -1.15%  "
$ws.Range("E29").Value = "  -0.34%  "
$ws.Range("D30").Formula = "=""9.12"""
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  +1.48%  "
$ws.Range("D31").Formula = "=""8.12"""
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = "  -0.93%  "
$ws.Range("E32").Value = "  -1.00%  "
$ws.Range("E33").Value = "  -0.15%  "
$ws.Range("D34").Formula = "=""6.81"""
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = "  -3.62%  "
$ws.Range("D35").Formula = "=""616.55"""
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  -8.22%  "
$ws.Range("D36").Formula = "=""0.0498"""
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = "  +5.12%  "
$ws.Range("E37").Value = "  +1.09%  "
$ws.Range("E38").Value = "  -0.78%  "
$ws.Range("D39").Formula = "=""56.89"""
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = "  -0.47%  "
$ws.Range("E40").Value = "  +0.21%  "
$ws.Range("D41").Formula = "=""3.39"""
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  -5.81%  "
$ws.Range("E42").Value = "  +1.99%  "
$ws.Range("D43").Value = "3.348.39"
$ws.Range("E43").Value = "  -0.72%  "
$ws.Range("D44").Value = "0.0₃0725"
$ws.Range("E44").Value = "  +3.90%  "
$ws.Range("E45").Value = "  -2.29%  "
$ws.Range("D46").Formula = "=""2.88"""
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = "  -1.63%  "
$ws.Range("D47").Formula = "=""31.86"""
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = "  -2.13%  "
$ws.Range("D48").Formula = "=""2.50"""
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = "  -3.12%  "
$ws.Range("E49").Value = "  -0.13%  "
$ws.Range("D50").Formula = "=""134.01"""
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = "  +1.41%  "
$ws.Range("D51").Formula = "=""0.155"""
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)

$excel.CutCopyMode = 0
